# Applies the "eventsliders" commit: adds four new Artisan command rows
# (showCurve / showExtraCurve / showEvents / showBackgroundEvents) to the
# "Commands" sheet, right before the existing "RC Command" section, and
# updates the two sheets' selections accordingly.

$wb = $excel.ActiveWorkbook

$sliders  = $wb.Worksheets.Item("Sliders")
$commands = $wb.Worksheets.Item("Commands")

# --- Commands sheet: insert 4 rows just above the old row 96 ("RC Command") ---
$commands.Activate() | Out-Null
$commands.Rows("96:99").Insert() | Out-Null

$commands.Range("B96").Value = "showCurve(<name>,<bool>)"
$commands.Range("C96").Value = "shows/hides the curve indicated by <name> which is one of { ET, BT, DeltaET, DeltaBT, BackgroundET, BackgroundBT}"

$commands.Range("B97").Value = "showExtraCurve(<extra_device>,<curve>,<bool>)"
$commands.Range("C97").Value = "shows/hides the <curve> (one of {T1,T2}) of the zero-based <extra_device> number"

$commands.Range("B98").Value = "showEvents(<event_type>, <bool>)"
$commands.Range("C98").Value = "shows/hides the events of <event_type> in [1,..,5]"

$commands.Range("B99").Value = "showBackgroundEvents(<bool>)"
$commands.Range("C99").Value = "shows/hides the events of the background profile"

# match the row height used by the rest of the "Commands" descriptive rows
$commands.Range("B96:C99").RowHeight = 13.8

# --- Selections, matching the post-edit workbook state ---
$sliders.Activate() | Out-Null
$sliders.Range("B6").Select() | Out-Null

$commands.Activate() | Out-Null
$commands.Range("C97").Select() | Out-Null
